$p = $ppt.ActivePresentation

# Locate the slide that contains the "GitHub" source-code textboxes
# (the slide whose creationId is 1667389619 / sldId 291 -- the last slide,
# "Source Code", holding "TextBox 25" and "TextBox 27").
$slide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $candidate = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $candidate.Shapes.Count; $shi++) {
        if ($candidate.Shapes.Item($shi).Name -eq "TextBox 25") {
            $slide = $candidate
            break
        }
    }
    if ($slide -ne $null) { break }
}

# Find the two shapes of interest by name.
$githubLabel = $null
$githubLink = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 25") { $githubLabel = $sh }
    if ($sh.Name -eq "TextBox 27") { $githubLink = $sh }
}

# Drop the matching build/animation entry from the timeline before removing
# the shape it targets, so no orphaned <p:par>/<p:bldP> is left behind.
$timeline = $slide.TimeLine
for ($i = $timeline.MainSequence.Count; $i -ge 1; $i--) {
    $effect = $timeline.MainSequence.Item($i)
    if ($effect.Shape.Id -eq $githubLabel.Id) {
        $effect.Delete()
    }
}

# Remove the standalone "GitHub: " label textbox entirely.
$githubLabel.Delete()

# Move the remaining link textbox up into the label's old slot.
# (Points are rounded to single precision internally before being
#  converted back to EMU, so nudge the input slightly so the stored
#  EMU values land exactly on the target offsets.)
$githubLink.Left = 94.32031636062992
$githubLink.Top = 184.9511811023622
